$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$words = @(
    'enden',
    'geben',
    'sichern',
    'wehtun',
    'stellen',
    'ahnen',
    'wachsen',
    'kratzen',
    'schlucken',
    'ändern',
    'wundern',
    'warnen',
    'spinnen',
    'wählen',
    'gelten',
    'achten',
    'dienen',
    'kümmern',
    'runden',
    'wirken',
    'husten',
    'trennen',
    'kaufen',
    'lesen',
    'jubeln',
    'segnen',
    'suchen',
    'klappen',
    'sinken',
    'spenden',
    'sorgen',
    'sperren',
    'tanzen',
    'schlafen',
    'meistern',
    'scheiden',
    'flüchten',
    'leihen',
    'ehren',
    'beten',
    'schwingen',
    'helfen',
    'pflanzen',
    'hassen',
    'zünden',
    'mögen',
    'treiben',
    'fließen',
    'führen',
    'zögern',
    'prügeln',
    'schreiten',
    'schultern',
    'sprechen',
    'äußern',
    'arten',
    'garen',
    'siegen',
    'freuen',
    'klingen',
    'stehlen',
    'greifen',
    'fehlen',
    'bleiben',
    'bluten',
    'dürfen',
    'pfeifen',
    'loben',
    'liegen',
    'rasen',
    'hacken',
    'jagen',
    'brauchen',
    'klingeln',
    'küssen',
    'schrecken',
    'fallen',
    'sagen',
    'hauen',
    'nähen',
    'lenken',
    'spüren',
    'lügen',
    'sparen',
    'tollen',
    'fischen',
    'münzen',
    'fällen',
    'biegen',
    'starren',
    'rücken',
    'werfen',
    'betteln',
    'lockern',
    'planen',
    'spielen',
    'steuern',
    'knarren',
    'regnen',
    'ärgern',
    'wagen',
    'heulen',
    'grüßen',
    'albern',
    'kichern',
    'hindern',
    'quälen',
    'duschen',
    'hoffen',
    'wollen',
    'irren',
    'fangen',
    'erben',
    'läuten',
    'altern',
    'scheitern',
    'rufen',
    'lächeln',
    'sitzen',
    'weichen',
    'folgen',
    'leugnen',
    'räumen',
    'bitten',
    'kehren',
    'platzen',
    'wüten',
    'wahren',
    'feiern',
    'malen',
    'hören',
    'heben',
    'schmecken',
    'reizen',
    'schaden',
    'trauen',
    'dauern',
    'drucken',
    'rechnen',
    'gründen',
    'drehen',
    'seufzen',
    'trotzen',
    'liefern',
    'heilen',
    'sterben',
    'mauern',
    'formen',
    'pflegen',
    'bellen',
    'boxen',
    'scheinen',
    'streichen',
    'bauen',
    'decken',
    'kennen',
    'werden',
    'handeln',
    'graben',
    'füllen',
    'zählen',
    'kosten',
    'dringen',
    'schenken',
    'stammen',
    'schicken',
    'machen',
    'backen',
    'löschen',
    'wenden',
    'klettern',
    'fahren',
    'stecken',
    'schwören',
    'filmen',
    'stechen',
    'kämpfen',
    'wecken',
    'schulden',
    'ruhen',
    'sprengen',
    'zeigen',
    'töten',
    'bergen',
    'tauschen',
    'zielen',
    'atmen',
    'tropfen',
    'flehen',
    'saufen',
    'schwächen',
    'pissen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
